$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.461.93'
$ws.Range('E2').Value = '  +0.89%  '

$ws.Range('D3').Value = '3.207.16'
$ws.Range('E3').Value = '  -0.89%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.32%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.05%  '

$ws.Range('E7').Value = '  -5.52%  '

$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('E9').Value = '  -2.44%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.73'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.91%  '

$ws.Range('E11').Value = '  -0.06%  '

$ws.Range('D12').Value = '3.763.38'
$ws.Range('E12').Value = '  -0.77%  '

$ws.Range('E13').Value = '  -0.42%  '

$ws.Range('D14').Value = '65.414.75'
$ws.Range('E14').Value = '  +0.75%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.71'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.40%  '

$ws.Range('D16').Value = '3.204.08'
$ws.Range('E16').Value = '  -0.06%  '

$ws.Range('E17').Value = '  -0.74%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '413.41'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.22%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.02%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.98%  '

$ws.Range('E21').Value = '  -0.65%  '

$ws.Range('E22').Value = '  +0.17%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.98%  '

$ws.Range('E24').Value = '  -1.80%  '

$ws.Range('E25').Value = '  -1.37%  '

$ws.Range('E26').Value = '  -5.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.20%  '

$ws.Range('E28').Value = '  +0.15%  '

$ws.Range('E29').Value = '  -1.40%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '21.65'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.08%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.03'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.34%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.41'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.58%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.14'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.68%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '157.46'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.05%  '

$ws.Range('E35').Value = '  -1.54%  '

$ws.Range('E36').Value = '  -0.04%  '

$ws.Range('D37').Value = '2.732.13'
$ws.Range('E37').Value = '  -3.23%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.29'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.58%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.16'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.93%  '

$ws.Range('E40').Value = '  -1.99%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0636'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.70%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.61'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.67%  '

$ws.Range('E43').Value = '  -0.53%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '297.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.04%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.65'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.56%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.00%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0990'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.34%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.99'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -8.37%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.82'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.14%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.44'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.48%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.911'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.34%  '
